# Applies the commit's edits:
#  1. Renumber the "STUDIES" section test items from 2.x to 1.x (A6:A17).
#  2. Add a new row content "Study-level Consent Details Report" (bold) in B19,
#     which previously was an empty placeholder row.
#  3. Remove the stray "2" value that had been left in A5.
#  4. Update the active selection to B19 (cursor moved down while editing).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Renumber test items 2.1 .. 2.12 -> 1.1 .. 1.12 (column A, rows 6-17) ---
for ($r = 6; $r -le 17; $r++) {
    $old = $ws.Cells.Item($r, 1).Value2
    if ($old -ne $null -and $old.ToString().StartsWith("2.")) {
        $new = $old.ToString() -replace '^2\.', '1.'
        $ws.Cells.Item($r, 1).Value = $new
    }
}

# --- 2. Clear the stray numeric "2" that used to sit in A5 ---
$ws.Cells.Item(5, 1).ClearContents()

# --- 3. Populate B19 with the new report-section heading, in bold ---
$ws.Cells.Item(19, 2).Value = "Study-level Consent Details Report"
$ws.Cells.Item(19, 2).Font.Bold = $true
$ws.Cells.Item(19, 2).HorizontalAlignment = -4131   # xlLeft
$ws.Cells.Item(19, 2).VerticalAlignment = -4160     # xlTop
$ws.Cells.Item(19, 2).WrapText = $true

# --- 4. Move the active selection down to B19 (matches the author's cursor) ---
$ws.Range("B19").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 13
